$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("D8").Value = "9 ماهه منتهی به 1399/09"
$ws.Range("E8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("F8").Value = "3 ماهه منتهی به 1400/03"
$ws.Range("G8").Value = "6 ماهه منتهی به 1400/06"
$ws.Range("H8").Value = "9 ماهه منتهی به 1400/09"
$ws.Range("I8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("J8").Value = "3 ماهه منتهی به 1401/03"
$ws.Range("K8").Value = "6 ماهه منتهی به 1401/06"
$ws.Range("L8").Value = "9 ماهه منتهی به 1401/09"
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"

$ws.Range("D9").Value = "1400-11-06 (3)"
$ws.Range("E9").Value = "1401-02-27 (11)"
$ws.Range("F9").Value = "1401-04-29 (3)"
$ws.Range("G9").Value = "1401-09-15 (5)"
$ws.Range("H9").Value = "1401-10-30 (3)"
$ws.Range("I9").Value = "1402-02-28 (7)"
$ws.Range("J9").Value = "1401-04-29"
$ws.Range("K9").Value = "1401-09-15 (2)"
$ws.Range("L9").Value = "1401-10-30"
$ws.Range("M9").Value = "1402-02-28"

$ws.Range("D11").Value = 12008
$ws.Range("E11").Value = 14557
$ws.Range("F11").Value = 2930
$ws.Range("G11").Value = 10609
$ws.Range("H11").Value = 15756
$ws.Range("I11").Value = 19296
$ws.Range("J11").Value = 3709
$ws.Range("K11").Value = 10524
$ws.Range("L11").Value = 17666
$ws.Range("M11").Value = 20800

$ws.Range("D12").Value = -7258
$ws.Range("E12").Value = -9416
$ws.Range("F12").Value = -2114
$ws.Range("G12").Value = -7095
$ws.Range("H12").Value = -10598
$ws.Range("I12").Value = -13316
$ws.Range("J12").Value = -2766
$ws.Range("K12").Value = -7414
$ws.Range("L12").Value = -12229
$ws.Range("M12").Value = -16725

$ws.Range("D13").Value = 4751
$ws.Range("E13").Value = 5140
$ws.Range("F13").Value = 816
$ws.Range("G13").Value = 3514
$ws.Range("H13").Value = 5158
$ws.Range("I13").Value = 5980
$ws.Range("J13").Value = 943
$ws.Range("K13").Value = 3110
$ws.Range("L13").Value = 5437
$ws.Range("M13").Value = 4075

$ws.Range("D14").Value = -792
$ws.Range("E14").Value = -1116
$ws.Range("F14").Value = -382
$ws.Range("G14").Value = -809
$ws.Range("H14").Value = -1198
$ws.Range("I14").Value = -1653
$ws.Range("J14").Value = -481
$ws.Range("K14").Value = -1025
$ws.Range("L14").Value = -1561
$ws.Range("M14").Value = -2009

$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "-"
$ws.Range("G15").Value = "-"
$ws.Range("H15").Value = "-"
$ws.Range("I15").Value = "-"
$ws.Range("J15").Value = "-"
$ws.Range("K15").Value = "-"
$ws.Range("L15").Value = "-"
$ws.Range("M15").Value = "-"

$ws.Range("D16").Value = 78
$ws.Range("E16").Value = -179
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -113
$ws.Range("I16").Value = -176
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = -9
$ws.Range("L16").Value = 34
$ws.Range("M16").Value = -294

$ws.Range("D17").Value = 4037
$ws.Range("E17").Value = 3846
$ws.Range("F17").Value = 443
$ws.Range("G17").Value = 2721
$ws.Range("H17").Value = 3847
$ws.Range("I17").Value = 4150
$ws.Range("J17").Value = 462
$ws.Range("K17").Value = 2076
$ws.Range("L17").Value = 3910
$ws.Range("M17").Value = 1773

$ws.Range("D18").Value = -962
$ws.Range("E18").Value = -1436
$ws.Range("F18").Value = -429
$ws.Range("G18").Value = -997
$ws.Range("H18").Value = -1491
$ws.Range("I18").Value = -1933
$ws.Range("J18").Value = -445
$ws.Range("K18").Value = -1168
$ws.Range("L18").Value = -1628
$ws.Range("M18").Value = -1965

$ws.Range("D19").Value = 594
$ws.Range("E19").Value = 5332
$ws.Range("F19").Value = 720
$ws.Range("G19").Value = 656
$ws.Range("H19").Value = 1540
$ws.Range("I19").Value = 7921
$ws.Range("J19").Value = 2340
$ws.Range("K19").Value = 4162
$ws.Range("L19").Value = 3951
$ws.Range("M19").Value = 11203

$ws.Range("D20").Value = 3669
$ws.Range("E20").Value = 7743
$ws.Range("F20").Value = 734
$ws.Range("G20").Value = 2381
$ws.Range("H20").Value = 3896
$ws.Range("I20").Value = 10138
$ws.Range("J20").Value = 2357
$ws.Range("K20").Value = 5070
$ws.Range("L20").Value = 6233
$ws.Range("M20").Value = 11010

$ws.Range("D21").Value = -803
$ws.Range("E21").Value = -653
$ws.Range("F21").Value = -3
$ws.Range("G21").Value = -467
$ws.Range("H21").Value = -796
$ws.Range("I21").Value = -536
$ws.Range("J21").Value = -3
$ws.Range("K21").Value = -189
$ws.Range("L21").Value = -463
$ws.Range("M21").Value = "-"

$ws.Range("D22").Value = 2866
$ws.Range("E22").Value = 7090
$ws.Range("F22").Value = 731
$ws.Range("G22").Value = 1914
$ws.Range("H22").Value = 3100
$ws.Range("I22").Value = 9602
$ws.Range("J22").Value = 2354
$ws.Range("K22").Value = 4882
$ws.Range("L22").Value = 5770
$ws.Range("M22").Value = 11010

$ws.Range("D23").Value = "-"
$ws.Range("E23").Value = "-"
$ws.Range("F23").Value = "-"
$ws.Range("G23").Value = "-"
$ws.Range("H23").Value = "-"
$ws.Range("I23").Value = 4
$ws.Range("J23").Value = "-"
$ws.Range("K23").Value = "-"
$ws.Range("L23").Value = "-"
$ws.Range("M23").Value = 5

$ws.Range("D24").Value = 2866
$ws.Range("E24").Value = 7090
$ws.Range("F24").Value = 731
$ws.Range("G24").Value = 1914
$ws.Range("H24").Value = 3100
$ws.Range("I24").Value = 9606
$ws.Range("J24").Value = 2354
$ws.Range("K24").Value = 4882
$ws.Range("L24").Value = 5770
$ws.Range("M24").Value = 11015

$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0

$ws.Range("D26").Value = 10311
$ws.Range("E26").Value = 10031
$ws.Range("F26").Value = 9721
$ws.Range("G26").Value = 9152
$ws.Range("H26").Value = 8718
$ws.Range("I26").Value = 8596
$ws.Range("J26").Value = 7718
$ws.Range("K26").Value = 7511
$ws.Range("L26").Value = 7132
$ws.Range("M26").Value = 6427

$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0
